# ---------------------------------------------------------------------------
# Applies the two changes captured in the commit diff:
#
#   1. Slide 5's table switches to a different built-in table style
#      ({8C073611-9E7F-4288-9D49-18983821D71F} -> {52380F04-E3F4-4BC2-893C-5689DC44882D}).
#
#   2. The deck's theme palette is swapped from the "Integral / Red Violet"
#      colors to the stock "Office Theme / Office" colors (the colors that
#      used to live only on the (inactive) notes-master theme become the
#      colors of the live design, mirroring the underlying theme1/theme2
#      part swap in the OOXML).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
$slide = $p.Slides.Item(5)

$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{52380F04-E3F4-4BC2-893C-5689DC44882D}")
}

# --- 2. Theme colors ------------------------------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
